# Auto-generated edit script applying the cryptos.xlsx symbol-list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.11%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "12.57%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.159"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.76%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07764"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.87%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.372"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6.01%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.024"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.48%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.947"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.37%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9324"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.87%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09899"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "8.24%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1795"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.61%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08604"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.13%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03315"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.06%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09921"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.53%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001503"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.13%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005768"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.18%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.461"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.18%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.06%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3367"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.20%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1334"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.18%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.295"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "7.62%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2300"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.41%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04565"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.08%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.28%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004366"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.34%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.17%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "10.70%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04800"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "6.10%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007746"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.60%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1414"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.51%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007137"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-27.52%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002099"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.03%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009178"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.20%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006113"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.40%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.10%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "34.19%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002000"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.10%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.10%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.10%"
